$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the "cmd window不該出現" row (old row 3). ---
# This shifts old rows 4-9 up to become rows 3-8, carrying their
# values/styles/row-heights with them (matches target rows 3,4,5,6 exactly
# with no further edits needed to those, except text/value tweaks below).
$ws.Rows(3).Delete()

# --- 2. Row 4 (was "revise .mp4 to .m4a(AAC or ALAC), support  FLAC(.flac)") ---
# New text; "done" column stays blank.
$ws.Range("A4").Value = "support  AAC(.m4a), FLAC(.flac), WMA(.wma), WAV(wav), MP3(.mp3)"

# --- 3. Row 5 (was "在加入queue前可選擇下載影片或音檔") ---
# New text; now marked done.
$ws.Range("A5").Value = "需要在convert時鎖定所有輸入"
$ws.Range("C5").Value = "v"

# --- 4. Row 6 (was "更改console成多個blocks，可刪除或修改block" / "use listbox") ---
# New text; clear the old remark; mark done.
$ws.Range("A6").Value = "禁止螢幕縮放"
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = "v"

# --- 5. Row 7 (was "ping 的取代方案") ---
# New text; row got visibly taller (21 vs 19.5).
$ws.Range("A7").Value = "允許使用playlist url下載整個playlist內的音檔"
$ws.Rows(7).RowHeight = 21

# --- 6. Row 8 (was "pytube再次出問題，目前等待修復中", unbordered style) ---
# Bring its formatting in line with the other bordered data rows before
# changing its text/remark, then fill in the new content.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("A8").Value = "執行效能顯然有待加強，集中管理未完成的urls、Youtube objects，只取得Youtube objects一次"
$ws.Range("B8").Value = "urls已完成"

# --- 7. Column widths: content got longer in col A & B, shorter in C.
# (ColumnWidth is quantized by the host to whole-pixel steps, so these
# inputs are chosen to land as close as possible to the real target
# widths of 105.58 / 12.43 / 6.72.)
$ws.Columns("A").ColumnWidth = 104.666
$ws.Columns("B").ColumnWidth = 11.6665
$ws.Columns("C").ColumnWidth = 5.834

$excel.CutCopyMode = 0
